$d = $word.ActiveDocument

# Locate the paragraph "Then edited the work document." — the four new
# paragraphs get inserted immediately after it, before the _GoBack bookmark.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq "Then edited the work document.") {
        $anchor = $candidate
        break
    }
}

$texts = @(
    "At the command line I did a:",
    "Git commit “A simple Word Document.docx”",
    "Then I edited the file again.",
    "Next I will commit those changes and try to do a  push back to github."
)

foreach ($t in $texts) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $t
}

# $anchor is now the last inserted paragraph ("Next I will commit...").
# The _GoBack bookmark paragraph immediately follows it.
$bookmarkPara = $anchor.Next()

# Insert one new, truly empty paragraph right after the bookmark paragraph
# (matches the extra "<w:p/>" the diff adds before the document's existing
# trailing blank paragraph).
$bookmarkPara.Range.InsertParagraphAfter()
$blankPara = $bookmarkPara.Next()
$blankPara.Range.Text = "TEMP"
$blankRange = $blankPara.Range
$trimmedRange = $d.Range($blankRange.Start, $blankRange.End - 1)
$trimmedRange.Delete()

Write-Output "done"
